# Remove duplicate CEDS sectors in Canada scaling mapping
#
# The "map" worksheet lists, for each inv_sector, a scaling_sector and one or
# more ceds_sector rows. A number of rows had the *same* ceds_sector value
# duplicated unnecessarily (sometimes even duplicated across multiple rows
# that map to the exact same scaling_sector). This script removes those
# duplicate ceds_sector (column C) entries:
#   - Most duplicates are simply cleared and given the existing "blank
#     placeholder" look (the grey-filled style already used elsewhere in the
#     sheet for intentionally empty ceds_sector cells, e.g. C34).
#   - The rows that exclusively used the now fully-redundant
#     "1A3eii_Other-transp" ceds_sector string are cleared back to a plain,
#     unstyled blank cell (removing the last references to that shared
#     string so it drops out of the workbook entirely).
#   - C73 keeps its existing style and is simply cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("map")

# --- Cells to clear completely (no style, no value) ---------------------
# These were the only remaining cells referencing the ceds_sector string
# "1A3eii_Other-transp"; clearing them drops that now-unused shared string.
$plainClear = @("C50", "C69", "C70", "C71")
foreach ($addr in $plainClear) {
    $ws.Range($addr).ClearContents()
}

# --- Cell that keeps its current style, only the value is removed -------
$ws.Range("C73").ClearContents()

# --- Cells that become blank "placeholder" cells -------------------------
# Reuse the grey-fill placeholder style already present on C34 (and others)
# instead of inventing a new one.
$styleSource = $ws.Range("C34")
$placeholderTargets = @("C23", "C25", "C26", "C27", "C52", "C53", "C54", "C58", "C59", "C62", "C64", "C65")
foreach ($addr in $placeholderTargets) {
    $target = $ws.Range($addr)
    $target.ClearContents()
    $styleSource.Copy()
    $target.PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# --- Restore the view/selection state recorded in the saved workbook ----
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 55
$ws.Range("C74").Select()
